# Trade #99 closed at 2026-02-17 15:57:42 - unknown UNKNOWN +0.000%
#
# - Summary sheet: Total Trades 98 -> 99, Win Rate % 36.73 -> 36.36
# - Strategy Status sheet (MarketMaking row): Trades 98 -> 99, Win Rate % 36.73 -> 36.36
# - All Trades sheet: append the newly closed trade as row 100
# - MarketMaking sheet: append the same newly closed trade as row 100

$wb = $excel.ActiveWorkbook

# --- Summary sheet ----------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B6").Value = 99      # Total Trades
$summary.Range("B9").Value = 36.36   # Win Rate %

# --- Strategy Status sheet ---------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D4").Value = 99       # MarketMaking Trades
$status.Range("G4").Value = 36.36    # MarketMaking Win Rate %

# --- New trade row data ------------------------------------------------------
$tradeNum      = 99
$tradeDate     = "2026-02-17"
$tradeTime     = "15:57:36"
$strategy      = "MarketMaking"
$side          = "DOWN"
$entryPrice    = 0.85
$exitPrice     = 0.85
$tradeStatus   = "CLOSED"
$pnlPct        = 0
$pnlDollar     = 0
$capitalAfter  = 99.90000000000001
$entrySlippage = 0
$exitSlippage  = 0
$confidence    = 0.6
$entryReason   = "Normal spread capture: 19600 bps"
$exitReason    = "early_exit"
$durationMin   = 0.14

function Add-TradeRow($ws, $row) {
    $ws.Cells.Item($row, 1).Value  = $tradeNum

    # Date/Time look like dates to the COM layer's auto-detection, so they
    # get typed in with a leading apostrophe to force plain text (the way
    # Excel itself avoids re-parsing a typed value as a date), then the
    # cell style is reset to Normal so no quote-prefix formatting lingers.
    $ws.Cells.Item($row, 2).Value = "'" + $tradeDate
    $ws.Cells.Item($row, 2).Style = "Normal"
    $ws.Cells.Item($row, 3).Value = "'" + $tradeTime
    $ws.Cells.Item($row, 3).Style = "Normal"

    $ws.Cells.Item($row, 4).Value  = $strategy
    $ws.Cells.Item($row, 5).Value  = $side
    $ws.Cells.Item($row, 6).Value  = $entryPrice
    $ws.Cells.Item($row, 7).Value  = $exitPrice
    $ws.Cells.Item($row, 8).Value  = $tradeStatus
    $ws.Cells.Item($row, 9).Value  = $pnlPct
    $ws.Cells.Item($row, 10).Value = $pnlDollar
    $ws.Cells.Item($row, 11).Value = $capitalAfter
    $ws.Cells.Item($row, 12).Value = $entrySlippage
    $ws.Cells.Item($row, 13).Value = $exitSlippage
    $ws.Cells.Item($row, 14).Value = $confidence
    $ws.Cells.Item($row, 15).Value = $entryReason
    $ws.Cells.Item($row, 16).Value = $exitReason
    $ws.Cells.Item($row, 17).Value = $durationMin
}

# --- All Trades sheet --------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $allTrades 100

# --- MarketMaking sheet -------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow $marketMaking 100
